# Add branch wise stocks
# Reorders the (Item Name, UOM) pairs within a few brand groups so the
# item listing comes out in the same order the branch-wise stock sheet uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Dinafex 180mg Tablet"
$ws.Range("E3").Value = "30's"
$ws.Range("D4").Value = "Dinafex 120mg Tablet"
$ws.Range("E4").Value = "30's"
$ws.Range("D5").Value = "Dinafex 60mg Tablet"
$ws.Range("E5").Value = "30's"

$ws.Range("D7").Value = "Etorix 120mg Tablet"
$ws.Range("E7").Value = "20's"
$ws.Range("D8").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E8").Value = "40's"
$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

$ws.Range("D11").Value = "Flucloxin 500mg Capsule"
$ws.Range("E11").Value = "30 's"
$ws.Range("D12").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E12").Value = "36 's"

$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

$ws.Range("D17").Value = "Kynol D 25mg Tablet"
$ws.Range("E17").Value = "60 's"
$ws.Range("D19").Value = "Kynol TR 200mg Capsule"
$ws.Range("E19").Value = "30 's"

$ws.Range("D24").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E24").Value = "30ml"
$ws.Range("D25").Value = "Zithrox 15ml Suspension"
$ws.Range("E25").Value = "15 ml"
$ws.Range("D27").Value = "Zithrox 500mg Tablet"
$ws.Range("E27").Value = "6 's"
